# chore: adapt column header formatting to respective input file names
#
# Renames the paired "_old"/"_new" header-column suffixes to the concrete
# format-version identifiers "_FV2310" / "_FV2404", turns the data range
# into a real Table ("Table1") and freezes the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename the "before" (FV2310) header columns A1:J1.
$fv2310Headers = @(
    "Segmentname_FV2310",
    "Segmentgruppe_FV2310",
    "Segment_FV2310",
    "Datenelement_FV2310",
    "Segment ID_FV2310",
    "Code_FV2310",
    "Qualifier_FV2310",
    "Beschreibung_FV2310",
    "Bedingungsausdruck_FV2310",
    "Bedingung_FV2310"
)
for ($i = 0; $i -lt $fv2310Headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $fv2310Headers[$i]
}

# Column K1 stays "diff" (unchanged).

# 2) Rename the "after" (FV2404) header columns L1:U1.
$fv2404Headers = @(
    "Segmentname_FV2404",
    "Segmentgruppe_FV2404",
    "Segment_FV2404",
    "Datenelement_FV2404",
    "Segment ID_FV2404",
    "Code_FV2404",
    "Qualifier_FV2404",
    "Beschreibung_FV2404",
    "Bedingungsausdruck_FV2404",
    "Bedingung_FV2404"
)
for ($i = 0; $i -lt $fv2404Headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = $fv2404Headers[$i]
}

# 3) Convert the used range A1:U79 into a proper table ("Table1"),
#    picking up the renamed header row.
$dataRange = $ws.Range("A1:U79")
$tbl = $ws.ListObjects.Add(
    [Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange,
    $dataRange,
    $null,
    [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes
)
$tbl.Name = "Table1"

# 4) Freeze the header row (split below row 1, top-left cell of the
#    scrollable pane is A2).
[void]$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
